# Auto-applies six-week O&M window + separate cable lay vessel updates
# to cumulative-capacity.xlsx: updates cumulative capacity values and
# removes the final (now out-of-horizon) row on the affected sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet: Baseline-Low ---
$ws = $wb.Worksheets.Item("Baseline-Low")

$values = @(
    1789.841500945225,
    2585.23032705802,
    3385.778190742603,
    4122.220710171389,
    4882.244586452645,
    5660.236667235732,
    6430.737169400287,
    7197.056682430292,
    8011.364786409796,
    9010.674642253591,
    10014.49205757578,
    11016.68207515606,
    11945.78359650867,
    12869.62190811221,
    13910.76788544547,
    15071.6429396868,
    16065.49263538764,
    17049.25545867653,
    18038.76765767459,
    19030.47256176306,
    20105.18084692169,
    21159.64732010838,
    22116.49740373248,
    23075.58070431934,
    24142.16066517226
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $values[$i]
}
$ws.Rows.Item(27).Delete()

# --- Sheet: Baseline-Mid (SC) ---
$ws = $wb.Worksheets.Item("Baseline-Mid (SC)")

$values = @(
    1020.793950850662,
    1728.285742712168,
    2657.520757058101,
    3750.343603486754,
    5427.658713295086,
    7309.043095529309,
    9198.492409208895,
    11102.82935368304,
    12962.37417889951,
    14874.12441014114,
    16834.06111285813,
    18613.94028419859,
    20451.55026848424,
    22311.37586101646,
    23382.21288438006,
    24502.3583666581
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $values[$i]
}

# --- Sheet: Baseline-Mid (CC) ---
$ws = $wb.Worksheets.Item("Baseline-Mid (CC)")

$values = @(
    1020.793950850662,
    1728.285742712168,
    2292.857854367965,
    3022.014144835143,
    3970.844176971362,
    5120.698435794835,
    6275.546487038782,
    7801.026650989696,
    9224.317300905135,
    10642.39031562669,
    12060.8830247035,
    13475.9649576049,
    14891.0718144499,
    16306.1786712949,
    17725.78312445612,
    19262.37738258394,
    20253.37883959044,
    20851.33105802048,
    21426.08359133127,
    21991.09907120743,
    22513.00053549116,
    23033.43023255814,
    23544.12790697674,
    24053.51423746579,
    24563.61730623653
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $values[$i]
}
$ws.Rows.Item(27).Delete()

# --- Sheet: Moderate-Low ---
$ws = $wb.Worksheets.Item("Moderate-Low")

$values = @(
    2100.788090051517,
    3360.693045196571,
    4621.066176147867,
    5814.115419579782,
    7025.568864549034,
    8245.456508266017,
    9462.154934163585,
    11293.72622110342,
    13184.73675914425,
    15316.19694713837,
    17493.03617798383,
    19750.51157410862,
    21445.47684804241,
    22869.62190811221,
    23910.76788544547,
    25071.64293968679,
    26065.49263538764,
    27049.25545867653,
    28038.76765767459,
    29030.47256176306,
    30105.18084692169,
    31159.64732010838,
    32116.49740373248,
    33075.58070431934,
    34142.16066517226
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $values[$i]
}
$ws.Rows.Item(27).Delete()

# --- Sheet: Moderate-Mid (SC) ---
$ws = $wb.Worksheets.Item("Moderate-Mid (SC)")

$values = @(
    1020.793950850662,
    2039.232331818461,
    3432.983475196653,
    4985.631588892019,
    6752.393779772964,
    8708.089550481423,
    10660.21609948145,
    12575.1801344746,
    15022.8644588079,
    17481.99944329268,
    19993.11046734116,
    22540.30641915152,
    25184.70384606915,
    27424.41408537074,
    29402.66073679991,
    31009.5540103765,
    31993.94343501695,
    32489.42068341254,
    32986.36569359584,
    33486.27204749129,
    33986.3005879331,
    34486.30136986302,
    34987.67123287671
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $values[$i]
}
$ws.Rows.Item(25).Delete()

# --- Sheet: Expanded-High ---
$ws = $wb.Worksheets.Item("Expanded-High")

$values = @(
    1020.793950850662,
    1728.285742712168,
    3476.842044841405,
    5386.747596299856,
    7519.306667389053,
    10294.99149913791,
    13171.92358050162,
    16228.91888841308,
    19647.12292875741,
    23442.09353861278,
    27634.67222762048,
    32200.56862013246,
    36705.95236495582,
    41219.71149603026,
    45680.72201222426,
    49634.36145522023,
    52918.53384826249,
    54632.89608898023
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $values[$i]
}
